$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2352941176470588
$ws.Range("C2").Value = 0.3529411764705883
$ws.Range("P2").Value = 0.2941176470588235
$ws.Range("S2").Value = 0.1176470588235294
# Row 3
$ws.Range("B3").Value = 0.1666666666666667
$ws.Range("J3").Value = 0.1666666666666667
$ws.Range("P3").Value = 0.6666666666666666
# Row 6
$ws.Range("J6").Value = 0.3846153846153846
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.07692307692307693
$ws.Range("S6").Value = 0.3846153846153846
# Row 7
$ws.Range("F7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("Q7").Value = 0.2222222222222222
$ws.Range("R7").Value = 0.2222222222222222
$ws.Range("S7").Value = 0.2777777777777778
# Row 8
$ws.Range("B8").Value = 0.04838709677419355
$ws.Range("D8").Value = 0.01612903225806452
$ws.Range("F8").Value = 0.01612903225806452
$ws.Range("J8").Value = 0.1612903225806452
$ws.Range("Q8").Value = 0.0967741935483871
$ws.Range("R8").Value = 0.0967741935483871
$ws.Range("S8").Value = 0.5645161290322581
# Row 9
$ws.Range("F9").Value = 0.1111111111111111
$ws.Range("J9").Value = 0.1851851851851852
$ws.Range("Q9").Value = 0.1851851851851852
$ws.Range("R9").Value = 0.03703703703703703
$ws.Range("S9").Value = 0.4814814814814815
# Row 10
$ws.Range("B10").Value = 0.0967741935483871
$ws.Range("D10").Value = 0.01075268817204301
$ws.Range("F10").Value = 0.03225806451612903
$ws.Range("J10").Value = 0.06451612903225806
$ws.Range("O10").Value = 0.01075268817204301
$ws.Range("Q10").Value = 0.1935483870967742
$ws.Range("R10").Value = 0.06451612903225806
$ws.Range("S10").Value = 0.5268817204301075
# Row 11
$ws.Range("G11").Value = 0.1363636363636364
$ws.Range("J11").Value = 0.04545454545454546
$ws.Range("K11").Value = 0.1818181818181818
$ws.Range("L11").Value = 0.5909090909090909
$ws.Range("S11").Value = 0.04545454545454546
# Row 12
$ws.Range("G12").Value = 0.8461538461538461
$ws.Range("J12").Value = 0.07692307692307693
$ws.Range("S12").Value = 0.07692307692307693
# Row 13
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.125
# Row 15
$ws.Range("H15").Value = 0.2307692307692308
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.3846153846153846
$ws.Range("K15").Value = 0.07692307692307693
$ws.Range("S15").Value = 0.2307692307692308
# Row 16
$ws.Range("H16").Value = 0.2727272727272727
$ws.Range("J16").Value = 0.2727272727272727
$ws.Range("K16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.3636363636363636
# Row 17
$ws.Range("F17").Value = 0.02857142857142857
$ws.Range("H17").Value = 0.2857142857142857
$ws.Range("I17").Value = 0.1142857142857143
$ws.Range("J17").Value = 0.2571428571428571
$ws.Range("K17").Value = 0.08571428571428572
$ws.Range("M17").Value = 0.05714285714285714
$ws.Range("O17").Value = 0.08571428571428572
$ws.Range("S17").Value = 0.08571428571428572
# Row 18
$ws.Range("H18").Value = 0.1764705882352941
$ws.Range("I18").Value = 0.1764705882352941
$ws.Range("J18").Value = 0.3529411764705883
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("O18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.1176470588235294
# Row 19
$ws.Range("F19").Value = 0.0145985401459854
$ws.Range("H19").Value = 0.3284671532846715
$ws.Range("I19").Value = 0.1313868613138686
$ws.Range("J19").Value = 0.2700729927007299
$ws.Range("K19").Value = 0.08029197080291971
$ws.Range("M19").Value = 0.04379562043795621
$ws.Range("O19").Value = 0.0364963503649635
$ws.Range("S19").Value = 0.0948905109489051
